$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

$ws.Range("A51").Value = 45139
$ws.Range("A52").Value = 45170
$ws.Range("A53").Value = 45200
$ws.Range("A54").Value = 45231
$ws.Range("A55").Value = 45261
$ws.Range("A56").Value = 45292
$ws.Range("A57").Value = 45323
$ws.Range("A58").Value = 45352
$ws.Range("A59").Value = 45383

$ws.Range("C49").Value = 1.25
$ws.Range("C50").Value = 1.25
$ws.Range("C51").Value = 1.25
$ws.Range("C52").Value = 1.25
$ws.Range("C53").Value = 1.25
$ws.Range("C54").Value = 1.25
$ws.Range("C55").Value = 1.25

$ws.Range("B55").Value = "FL(5-0-0)"
$ws.Range("D55").Value = 5

$ws.Activate()
$ws.Range("E55").Select()
